$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E10").Value = "both"
$ws.Range("E2:E10").Select()
